$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'43.519.71"
$ws.Range("E2").Value = "  +2.38%  "
$ws.Range("D3").Value = "'2.252.68"
$ws.Range("E3").Value = "  +2.28%  "
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").Value = "'321.00"
$ws.Range("E5").Value = "  +3.32%  "
$ws.Range("D6").Value = "'101.69"
$ws.Range("E6").Value = "  +4.11%  "
$ws.Range("D7").Value = "'0.585"
$ws.Range("E7").Value = "  +2.09%  "
$ws.Range("E8").Value = "  -0.23%  "
$ws.Range("D9").Value = "'0.567"
$ws.Range("E9").Value = "  +2.70%  "
$ws.Range("D10").Value = "'37.75"
$ws.Range("E10").Value = "  +4.03%  "
$ws.Range("D11").Value = "'0.0843"
$ws.Range("E11").Value = "  +2.51%  "
$ws.Range("D12").Value = "'7.80"
$ws.Range("E12").Value = "  +3.63%  "
$ws.Range("E13").Value = "  +3.00%  "
$ws.Range("D14").Value = "'0.872"
$ws.Range("E14").Value = "  +2.91%  "
$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").Value = "'2.588.71"
$ws.Range("E15").Value = "  +1.91%  "
$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D16").Value = "'14.44"
$ws.Range("E16").Value = "  +4.47%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "'2.252.73"
$ws.Range("E17").Value = "  +2.71%  "
$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").Value = "'43.451.03"
$ws.Range("E18").Value = "  +2.56%  "
$ws.Range("B19").Value = "InternetComputer(DFINITY)"
$ws.Range("C19").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D19").Value = "'14.21"
$ws.Range("E19").Value = "  +0.26%  "
$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").Value = "'0.0000101"
$ws.Range("E20").Value = "  +7.55%  "
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").Value = "'6.65"
$ws.Range("E21").Value = "  +2.99%  "
$ws.Range("B22").Value = "Litecoin"
$ws.Range("C22").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D22").Value = "'65.87"
$ws.Range("E22").Value = "  +1.96%  "
$ws.Range("B23").Value = "PancakeSwap"
$ws.Range("C23").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D23").Value = "'3.18"
$ws.Range("E23").Value = "  +0.57%  "
$ws.Range("B24").Value = "BitcoinCash"
$ws.Range("C24").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D24").Value = "'237.52"
$ws.Range("E24").Value = "  +2.40%  "
$ws.Range("B25").Value = "ImmutableX"
$ws.Range("C25").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D25").Value = "'2.21"
$ws.Range("E25").Value = "  +6.43%  "
$ws.Range("B26").Value = "Dai"
$ws.Range("C26").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  +0.13%  "
$ws.Range("B27").Value = "LEO"
$ws.Range("C27").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D27").Value = "'4.06"
$ws.Range("E27").Value = "  +2.92%  "
$ws.Range("B28").Value = "Cosmos"
$ws.Range("C28").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D28").Value = "'10.17"
$ws.Range("E28").Value = "  +1.28%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").Value = "'2.22"
$ws.Range("E29").Value = "  +3.44%  "
$ws.Range("D30").Value = "'37.86"
$ws.Range("E30").Value = "  +13.70%  "
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").Value = "'6.44"
$ws.Range("E31").Value = "  +0.30%  "
$ws.Range("D32").Value = "'0.0882"
$ws.Range("E32").Value = "  +2.47%  "
$ws.Range("B33").Value = "Monero"
$ws.Range("C33").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D33").Value = "'161.28"
$ws.Range("E33").Value = "  +3.05%  "
$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D34").Value = "'20.42"
$ws.Range("E34").Value = "  +0.61%  "
$ws.Range("B35").Value = "WEMIXToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D35").Value = "'2.73"
$ws.Range("E35").Value = "  +1.41%  "
$ws.Range("B36").Value = "LidoDAOToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D36").Value = "'3.27"
$ws.Range("E36").Value = "  +3.64%  "
$ws.Range("D37").Value = "'1.95"
$ws.Range("E37").Value = "  +7.06%  "
$ws.Range("B38").Value = "Stellar"
$ws.Range("C38").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D38").Value = "'0.121"
$ws.Range("E38").Value = "  +0.24%  "
$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D39").Value = "'4.44"
$ws.Range("E39").Value = "  +0.56%  "
$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").Value = "'0.106"
$ws.Range("E40").Value = "  +2.56%  "
$ws.Range("B41").Value = "NEARProtocol"
$ws.Range("C41").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D41").Value = "'3.80"
$ws.Range("E41").Value = "  +8.87%  "
$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").Value = "'0.0324"
$ws.Range("E42").Value = "  +2.77%  "
$ws.Range("B43").Value = "Celestia"
$ws.Range("C43").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D43").Value = "'15.30"
$ws.Range("E43").Value = "  +29.58%  "
$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D44").Value = "'1.00"
$ws.Range("E44").Value = "  -0.09%  "
$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").Value = "'1.809.53"
$ws.Range("E45").Value = "  +2.25%  "
$ws.Range("B46").Value = "Algorand"
$ws.Range("C46").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D46").Value = "'0.206"
$ws.Range("E46").Value = "  +1.48%  "
$ws.Range("B47").Value = "BitcoinSV"
$ws.Range("C47").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D47").Value = "'84.73"
$ws.Range("E47").Value = "  -4.49%  "
$ws.Range("B48").Value = "THORChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D48").Value = "'5.34"
$ws.Range("E48").Value = "  +1.69%  "
$ws.Range("B49").Value = "ordi"
$ws.Range("C49").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D49").Value = "'75.50"
$ws.Range("E49").Value = "  -1.09%  "
$ws.Range("B50").Value = "FraxShare"
$ws.Range("C50").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D50").Value = "'8.84"
$ws.Range("E50").Value = "  +5.10%  "
$ws.Range("B51").Value = "MultiversX"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D51").Value = "'59.14"
$ws.Range("E51").Value = "  -0.70%  "
